# Update reaction_sensitivity values on both sheets (NBR and BAR)
$wb = $excel.ActiveWorkbook

$nbr = $wb.Worksheets.Item("NBR")
$bar = $wb.Worksheets.Item("BAR")

$nbrValues = @(872, 866, 869, 845, 845, 834, 832, 825, 835, 831, 828, 823, 821, 820, 797, 795, 797, 795, 794)
$barValues = @(711, 710, 707, 732, 726, 742, 729, 730, 713, 713, 709, 712, 709, 711, 706, 700, 696, 694, 694)

for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $row = $i + 2
    $nbr.Cells.Item($row, 3).Value = $nbrValues[$i]
}

for ($i = 0; $i -lt $barValues.Length; $i++) {
    $row = $i + 2
    $bar.Cells.Item($row, 3).Value = $barValues[$i]
}
